$d = $word.ActiveDocument

# --- Step 1: remove the "Meta description" paragraph (2nd paragraph) ---
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# --- Step 2: insert a new bold paragraph before the last paragraph ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newRange = $lastPara.Range.InsertParagraphBefore()

# The newly inserted (empty) paragraph is now the second-to-last paragraph.
$titlePara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$titleRange = $titlePara.Range
$titleRange.Text = "Play Book of Gods for Free - Review by Slot Expert"
$titleRange.Font.Bold = 1

# --- Step 3: replace the text of the last paragraph (Maya warrior prompt) ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
# Trim off the trailing paragraph mark so only the run text is replaced
$lastRange.MoveEnd(1, -1) | Out-Null
$lastRange.Text = "Read our review of Book of Gods - an Ancient Egyptian-themed online slot game. Play this visually amazing slot for free and explore its exciting features."
$lastRange.Font.Italic = 1
